# ECP-TABLE.xlsx edit: replace the ";"-separated validity/invalidity criteria
# with proper "&&" (AND, for the valid-ECP column) / "||" (OR, for the
# invalid-ECP column) boolean expressions, widen columns C/E to fit the
# longer text, move the selection, and append a new empty "Folha2" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- registerUser (rows 2-4) ---
$ws.Range("C3").Value = "IDUser = INT &&  name= String  && rentalProgram = INT"
$ws.Range("E3").Value = "IDUser!= INT || name!=String || rentalProgram!=INT"
$ws.Range("C4").Value = "IDUser >= 0 && name = ""name"" && rentalProgram = 1 && rentalProgram = 2"
$ws.Range("E4").Value = "IDUser<0 || name=null || rentalProgram !=1 || rentalProgram !=2"

# --- addCredit (rows 7-9) ---
$ws.Range("C8").Value = "IDUser = INT && amount = INT"
$ws.Range("E8").Value = "IDUser != INT || amount != INT"
$ws.Range("C9").Value = "IDUser>=0 && amount > 0"
$ws.Range("E9").Value = "IDUser<0 || amount<= 0"

# --- bicycleRentalFee (rows 17-19) ---
$ws.Range("C18").Value = "rentalProgram = INT && startTime = INT && endTime = INT && nRentals = INT"
$ws.Range("E18").Value = "rentalProgram != INT ||startTime != INT || endTime != INT || nRentals != INT"
$ws.Range("C19").Value = "rentalProgram>=0 && startTime>=0 && endTime>=0 && nRentals>=0 && startTime<=endTime"
$ws.Range("E19").Value = "rentalProgram<0 || startTime<0 || endTime<0 || nRentals<0 || startTime>endTime"

# --- returnBicycle (rows 22-24) ---
$ws.Range("C23").Value = "IDUser = INT && IDDeposit = INT && endTime = INT"
$ws.Range("E23").Value = "IDUser != INT || IDDeposit != INT || endTime != INT"
$ws.Range("C24").Value = "IDUser>=0 && IDDeposit>0 && endTime>=0"
$ws.Range("E24").Value = "IDUser<0 || IDDeposit<=0 || endTime<0"

# --- getBicycle (rows 27-29) ---
$ws.Range("C28").Value = "IDUser = INT && IDDeposit = INT && startTime = INT"
$ws.Range("E28").Value = "IDUser != INT || IDDeposit != INT || startTime != INT"
$ws.Range("C29").Value = "IDUser>0 && IDDeposit>0 && startTime>=0"
$ws.Range("E29").Value = "IDUser<=0 || IDDeposit<=0 || startTime<0"

# Widen columns C and E to fit the now-longer criteria text (closest the
# engine's column-width quantization can reach to the authored 79.332 / 70.441
# character widths).
$ws.Columns.Item(3).ColumnWidth = 78.45
$ws.Columns.Item(5).ColumnWidth = 69.65

# Add a second, blank worksheet ("Folha2") after "Folha1".
$newSheet = $wb.Worksheets.Add($null, $ws)
$newSheet.Name = "Folha2"

# Keep "Folha1" the active/selected sheet and tab, and move the selection/
# active cell like the author did before saving.
$ws.Activate()
$ws.Range("E26").Select()

Write-Output "done"
